$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.903.24'
$ws.Range("E2").Value = '  -0.53%  '
$ws.Range("D3").Value = '3.034.00'
$ws.Range("E3").Value = '  -0.58%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''586.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.21%  '
$ws.Range("D6").Value = '''149.44'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.92%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '''0.526'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.37%  '
$ws.Range("D9").Value = '3.035.65'
$ws.Range("E9").Value = '  -0.54%  '
$ws.Range("D10").Value = '''0.152'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.86%  '
$ws.Range("D11").Value = '''5.73'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.77%  '
$ws.Range("D12").Value = '''0.445'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.15%  '
$ws.Range("D13").Value = '''0.0000232'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.79%  '
$ws.Range("D14").Value = '''35.38'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.76%  '
$ws.Range("E15").Value = '  +2.28%  '
$ws.Range("D16").Value = '3.533.75'
$ws.Range("E16").Value = '  -0.67%  '
$ws.Range("E17").Value = '  -0.03%  '
$ws.Range("D18").Value = '62.848.44'
$ws.Range("E18").Value = '  -0.75%  '
$ws.Range("D19").Value = '3.030.16'
$ws.Range("E19").Value = '  -0.57%  '
$ws.Range("D20").Value = '''468.92'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.47%  '
$ws.Range("D21").Value = '''14.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.82%  '
$ws.Range("D22").Value = '''0.694'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.34%  '
$ws.Range("D23").Value = '''7.42'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.52%  '
$ws.Range("D24").Value = '''2.37'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.63%  '
$ws.Range("D25").Value = '''81.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.71%  '
$ws.Range("D26").Value = '''12.43'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.63%  '
$ws.Range("D27").Value = '''10.41'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.93%  '
$ws.Range("D28").Value = '''1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("E29").Value = '  -2.72%  '
$ws.Range("D30").Value = '''0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("D31").Value = '''2.64'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.42%  '
$ws.Range("D32").Value = '''2.15'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.33%  '
$ws.Range("D33").Value = '''27.56'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.99%  '
$ws.Range("E34").Value = '  -3.99%  '
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("D36").Value = '0.0₃0803'
$ws.Range("E36").Value = '  -2.16%  '
$ws.Range("D37").Value = '''5.80'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.48%  '
$ws.Range("D38").Value = '''2.16'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.33%  '
$ws.Range("D39").Value = '''50.34'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.74%  '
$ws.Range("E40").Value = '  -2.78%  '
$ws.Range("D41").Value = '''2.97'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -10.10%  '
$ws.Range("D42").Value = '''425.86'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.82%  '
$ws.Range("E43").Value = '  +2.62%  '
$ws.Range("E44").Value = '  -0.93%  '
$ws.Range("D45").Value = '2.805.71'
$ws.Range("E45").Value = '  +0.73%  '
$ws.Range("E46").Value = '  -0.12%  '
$ws.Range("D47").Value = '''37.85'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.85%  '
$ws.Range("D48").Value = '''129.01'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.03%  '
$ws.Range("D50").Value = '''24.53'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.34%  '
$ws.Range("E51").Value = '  -0.29%  '
